$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has two Pearson logo pictures living in the page footers
# (both currently named "image1.png") and one BTEC logo picture living in
# the "first page" header (currently named "image2.jpg"). Renaming the
# footer inline shapes directly (InlineShape.Name = ...) hits a stale
# "addressed block" resolution issue for footer stories in this engine,
# so we route those two renames through the Selection object instead --
# selecting the shape's range first, then renaming via
# $word.Selection.InlineShapes -- which resolves the block correctly.

for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $footer = $sec.Footers.Item($f)
    if ($footer.Exists -and $footer.Range.InlineShapes.Count -gt 0) {
        $pic = $footer.Range.InlineShapes.Item(1)
        $pic.Range.Select()
        $word.Selection.InlineShapes.Item(1).Name = "image2.png"
    }
}

for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $header = $sec.Headers.Item($h)
    if ($header.Exists -and $header.Range.InlineShapes.Count -gt 0) {
        $pic = $header.Range.InlineShapes.Item(1)
        $pic.Name = "image1.jpg"
    }
}
